$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare helper cells (off to the side) carrying the exact styles we need,
# so that PasteSpecial(xlPasteFormats) can stamp them onto target cells without
# Excel synthesizing brand-new style records.
$ws.Range("AA1").Font.Bold = $true          # matches existing bold/no-alignment style
# AA2 is left completely untouched -> default style (0)

# ===================== Row 1 header: F1 "Bump Sensor" -> "Voltage" =====================
$ws.Range("AA1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F1").Value = "Voltage"

# ===================== Rows 22-23: Tiva Port C pins, F column "Left"/"Right" -> "3.3 V" =====================
$ws.Range("A22").Value = "C"

$ws.Range("AA2").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (default style, no bold)
$ws.Range("F22").Value = "3.3 V"

$ws.Range("AA2").Copy() | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Value = "3.3 V"

# ===================== New rows 25-27: Tiva Port D pins 0,1,2 wired to 5V =====================
$ws.Range("A25").Value = "D"
$ws.Range("B25").Value = 0
$ws.Range("AA2").Copy() | Out-Null
$ws.Range("F25").PasteSpecial(-4122) | Out-Null
$ws.Range("F25").Value = "5V"

$ws.Range("B26").Value = 1
$ws.Range("AA2").Copy() | Out-Null
$ws.Range("F26").PasteSpecial(-4122) | Out-Null
$ws.Range("F26").Value = "5V"

$ws.Range("B27").Value = 2
$ws.Range("AA2").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = "5V"

# --- cleanup helper cells
$ws.Range("AA1:AA2").Clear() | Out-Null

# ===================== Column widths: F narrower & loses its style, G loses bestFit =====================
$ws.Columns("F").ColumnWidth = 8.0533854166667
$ws.Columns("G").ColumnWidth = 43.75

# ===================== View state: scroll so row 4 is at top, select F27 =====================
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("F27").Select() | Out-Null
